# Added pricing tab test cases
# - Inserts a new "budgetTest" row before the existing "chainSummaryTest"
#   row (row 10) on both the "Sheet1" and "Data" worksheets, pushing the
#   old row 10 down to row 11.
# - Flips the Execute flag for a couple of existing test cases from
#   "yes" to "no" (createOccupancyRule on Sheet1, createPromotionRule and
#   createOccupancyRule on Data) now that the new budget test case is the
#   one that executes.
# - Leaves the selection on the newly-added pricing test rows, with the
#   "Sheet1" tab active (it had been on "Data").

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet1: TestCaseName | Test Case Description | Execute | InvocationCount | Priority
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Sheet1")

# Insert a new row 10, pushing the existing row 10 ("chainSummaryTest")
# down to row 11.
$ws1.Rows.Item(10).Insert()

# Match the formatting of the row directly above (row 9) for the new row.
$ws1.Range("B9:E9").Copy()
$ws1.Range("B10:E10").PasteSpecial(-4122)

$ws1.Range("A10").Value = "budgetTest"
$ws1.Range("B10").Value = "abcd"
$ws1.Range("C10").Value = "yes"
$ws1.Range("D10").Value = "'1"
$ws1.Range("E10").Value = "'1"

# createOccupancyRule no longer executes.
$ws1.Range("C9").Value = "no"

# ---------------------------------------------------------------------
# Data: TestName | Execute | Browser | UserName | Password
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Data")

# Insert a new row 10, pushing the existing row 10 ("chainSummaryTest")
# down to row 11.
$ws2.Rows.Item(10).Insert()

# Match the formatting of the row directly above (row 9) for the new row.
$ws2.Range("B9:E9").Copy()
$ws2.Range("B10:E10").PasteSpecial(-4122)

$ws2.Range("A10").Value = "budgetTest"
$ws2.Range("B10").Value = "yes"
$ws2.Range("C10").Value = "chrome"
$ws2.Range("D10").Value = "raghavendra.m@axisrooms.com"
$ws2.Range("E10").Value = "Password123#"

# createPromotionRule / createOccupancyRule no longer execute.
$ws2.Range("B8").Value = "no"
$ws2.Range("B9").Value = "no"

# ---------------------------------------------------------------------
# Selection / active tab bookkeeping to mirror the authored edit.
# ---------------------------------------------------------------------
$ws2.Activate()
$ws2.Range("A10").Select()

$ws1.Activate()
$ws1.Range("A16").Select()
